$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 228
$ws.Range("A228").Value = "'2025-07-11"
$ws.Range("B228").Value = "Llaneros"
$ws.Range("C228").Value = "America de Cali"
$ws.Range("D228").Value = 0
$ws.Range("E228").Value = 0
$ws.Range("F228").Value = 1392345
$ws.Range("G228").Value = 6
$ws.Range("H228").Value = 7
$ws.Range("I228").Value = 3
$ws.Range("J228").Value = 1
$ws.Range("K228").Value = 0
$ws.Range("L228").Value = 0
$ws.Range("M228").Value = 0
$ws.Range("N228").Value = 0
$ws.Range("O228").Value = 0
$ws.Range("P228").Value = 0
$ws.Range("Q228").Value = 41
$ws.Range("R228").Value = 59
$ws.Range("S228").Value = "E"

# Row 229
$ws.Range("A229").Value = "'2025-07-12"
$ws.Range("B229").Value = "Once Caldas"
$ws.Range("C229").Value = "Atletico Nacional"
$ws.Range("D229").Value = 1
$ws.Range("E229").Value = 3
$ws.Range("F229").Value = 1392346
$ws.Range("G229").Value = 2
$ws.Range("H229").Value = 7
$ws.Range("I229").Value = 5
$ws.Range("J229").Value = 2
$ws.Range("K229").Value = 0
$ws.Range("L229").Value = 0
$ws.Range("M229").Value = 0
$ws.Range("N229").Value = 0
$ws.Range("O229").Value = 1
$ws.Range("P229").Value = 3
$ws.Range("Q229").Value = 39
$ws.Range("R229").Value = 61
$ws.Range("S229").Value = "V"

# Row 230
$ws.Range("A230").Value = "'2025-07-12"
$ws.Range("B230").Value = "Deportivo Pasto"
$ws.Range("C230").Value = "Deportes Tolima"
$ws.Range("D230").Value = 3
$ws.Range("E230").Value = 2
$ws.Range("F230").Value = 1392347
$ws.Range("G230").Value = 0
$ws.Range("H230").Value = 3
$ws.Range("I230").Value = 5
$ws.Range("J230").Value = 1
$ws.Range("K230").Value = 0
$ws.Range("L230").Value = 0
$ws.Range("M230").Value = 0
$ws.Range("N230").Value = 0
$ws.Range("O230").Value = 3
$ws.Range("P230").Value = 2
$ws.Range("Q230").Value = 23
$ws.Range("R230").Value = 77
$ws.Range("S230").Value = "L"

# Row 231
$ws.Range("A231").Value = "'2025-07-12"
$ws.Range("B231").Value = "Envigado"
$ws.Range("C231").Value = "Fortaleza FC"
$ws.Range("D231").Value = 0
$ws.Range("E231").Value = 1
$ws.Range("F231").Value = 1392348
$ws.Range("G231").Value = 5
$ws.Range("H231").Value = 1
$ws.Range("I231").Value = 3
$ws.Range("J231").Value = 2
$ws.Range("K231").Value = 0
$ws.Range("L231").Value = 0
$ws.Range("M231").Value = 0
$ws.Range("N231").Value = 0
$ws.Range("O231").Value = 0
$ws.Range("P231").Value = 1
$ws.Range("Q231").Value = 46
$ws.Range("R231").Value = 54
$ws.Range("S231").Value = "V"

# Row 232
$ws.Range("A232").Value = "'2025-07-13"
$ws.Range("B232").Value = "Bucaramanga"
$ws.Range("C232").Value = "Chico"
$ws.Range("D232").Value = 1
$ws.Range("E232").Value = 1
$ws.Range("F232").Value = 1392349
$ws.Range("G232").Value = 6
$ws.Range("H232").Value = 4
$ws.Range("I232").Value = 2
$ws.Range("J232").Value = 2
$ws.Range("K232").Value = 0
$ws.Range("L232").Value = 0
$ws.Range("M232").Value = 0
$ws.Range("N232").Value = 0
$ws.Range("O232").Value = 1
$ws.Range("P232").Value = 1
$ws.Range("Q232").Value = 68
$ws.Range("R232").Value = 32
$ws.Range("S232").Value = "E"

# Row 233
$ws.Range("A233").Value = "'2025-07-13"
$ws.Range("B233").Value = "Deportivo Cali"
$ws.Range("C233").Value = "Junior"
$ws.Range("D233").Value = 0
$ws.Range("E233").Value = 2
$ws.Range("F233").Value = 1392350
$ws.Range("G233").Value = 6
$ws.Range("H233").Value = 2
$ws.Range("I233").Value = 1
$ws.Range("J233").Value = 3
$ws.Range("K233").Value = 1
$ws.Range("L233").Value = 0
$ws.Range("M233").Value = 0
$ws.Range("N233").Value = 0
$ws.Range("O233").Value = 0
$ws.Range("P233").Value = 2
$ws.Range("Q233").Value = 27
$ws.Range("R233").Value = 73
$ws.Range("S233").Value = "V"

# Row 234
$ws.Range("A234").Value = "'2025-07-13"
$ws.Range("B234").Value = "Independiente Medellin"
$ws.Range("C234").Value = "Alianza Petrolera"
$ws.Range("D234").Value = 1
$ws.Range("E234").Value = 1
$ws.Range("F234").Value = 1392351
$ws.Range("G234").Value = 15
$ws.Range("H234").Value = 3
$ws.Range("I234").Value = 1
$ws.Range("J234").Value = 2
$ws.Range("K234").Value = 0
$ws.Range("L234").Value = 0
$ws.Range("M234").Value = 0
$ws.Range("N234").Value = 0
$ws.Range("O234").Value = 1
$ws.Range("P234").Value = 1
$ws.Range("Q234").Value = 67
$ws.Range("R234").Value = 33
$ws.Range("S234").Value = "E"

# Row 235
$ws.Range("A235").Value = "'2025-07-13"
$ws.Range("B235").Value = "Deportivo Pereira"
$ws.Range("C235").Value = "Santa Fe"
$ws.Range("D235").Value = 2
$ws.Range("E235").Value = 2
$ws.Range("F235").Value = 1392352
$ws.Range("G235").Value = 7
$ws.Range("H235").Value = 3
$ws.Range("I235").Value = 3
$ws.Range("J235").Value = 2
$ws.Range("K235").Value = 0
$ws.Range("L235").Value = 1
$ws.Range("M235").Value = 0
$ws.Range("N235").Value = 0
$ws.Range("O235").Value = 2
$ws.Range("P235").Value = 2
$ws.Range("Q235").Value = 69
$ws.Range("R235").Value = 31
$ws.Range("S235").Value = "E"

# Row 236
$ws.Range("A236").Value = "'2025-07-14"
$ws.Range("B236").Value = "La Equidad"
$ws.Range("C236").Value = "Rionegro Aguilas"
$ws.Range("D236").Value = 0
$ws.Range("E236").Value = 0
$ws.Range("F236").Value = 1392353
$ws.Range("G236").Value = 2
$ws.Range("H236").Value = 7
$ws.Range("I236").Value = 2
$ws.Range("J236").Value = 2
$ws.Range("K236").Value = 1
$ws.Range("L236").Value = 1
$ws.Range("M236").Value = 0
$ws.Range("N236").Value = 0
$ws.Range("O236").Value = 0
$ws.Range("P236").Value = 0
$ws.Range("Q236").Value = 46
$ws.Range("R236").Value = 54
$ws.Range("S236").Value = "E"
